$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for column D (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.078.22"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.874.44"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "312.93"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5052"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "0.08399"
$ws.Range("E9").Value = "  -8.29%  "
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "41.64"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "6.224"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "1.879.82"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "7.197"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "0.06656"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "6.062"
$ws.Range("D23").Value = "28.116.98"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "2.578"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").Value = "2.099.05"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "157.29"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.63"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "126.54"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "0.1050"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "1.050"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").Value = "5.618"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "9.702"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").Value = "0.02455"
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("D37").Value = "0.06549"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "0.2168"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "1.219"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "0.6507"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("E41").Value = "  -7.08%  "
$ws.Range("D42").Value = "11.33"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "4.894"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "0.6194"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").Value = "13.06"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "1.303"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "3.681"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "1.218"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "121.00"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "80.56"
$ws.Range("E51").Value = "  +2.32%  "

# Restore default style on column D so no stray explicit format/style remains
$ws.Range("D2:D51").Style = "Normal"

